# Atualizacao rapida de agenda as  8:38:50,98
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Giovani / 0893 / Auto Posto Aliança / ...
$ws.Range("A2").Value = "Giovani"
$ws.Range("B2").Value = "'0893"
$ws.Range("C2").Value = "Auto Posto Aliança"
$ws.Range("D2").Value = "Subir o monitoramento de câmeras pra base, possível reset de DVR."
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "Em andamento"
$ws.Range("I2").Value = "Fonte 10AP, módulo ethernet"

# Row 3: Roberto / 0773 / Escola Antônio Gonçalves de Matos / ...
$ws.Range("A3").Value = "Roberto"
$ws.Range("B3").Value = "'0773"
$ws.Range("C3").Value = "Escola Antônio Gonçalves de Matos"
$ws.Range("D3").Value = "Algumas câmeras estão sem imagem."
$ws.Range("E3").Value = "Roberto disse que foi ao local e o colégio estava fechado."
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "Falha"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""

# Row 2's H cell (reuses index allocated after row 3's new strings)
$ws.Range("H2").Value = "Maxvel: 20 / Forte: 14"

# Row 4: Roberto / 0706 / Lar das Meninas / ... (row no longer uses a tall custom height)
$ws.Range("A4").Value = "Roberto"
$ws.Range("B4").Value = "'0706"
$ws.Range("C4").Value = "Lar das Meninas"
$ws.Range("D4").Value = "Algumas câmeras estão sem imagem."
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "Em andamento"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Rows.Item(4).AutoFit()

# Rows 5-15: clear all the old content (row 8's I column keeps its note)
$ws.Range("A5:I15").ClearContents()
$ws.Range("I8").Value = "1 par de ballun IntelBras."

# Move the active selection to H2 (was H4)
$ws.Range("H2").Select()
